# Update fake SSN-like values so they all start with "9" (to clearly mark
# them as fake data), as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new value (first digit/character replaced with "9")
$ws.Range("B1").Value = 946412419
$ws.Range("B3").Value = 909360961
$ws.Range("B4").Value = 987777434
$ws.Range("B5").Value = 973351423
$ws.Range("B6").Value = 909175891
$ws.Range("B7").Value = 909175891
$ws.Range("B10").Value = 965507
$ws.Range("B11").Value = 969005507
$ws.Range("B12").Value = 969860000

# The shared string "44641241A" stored in B13 becomes "94641241A"
$ws.Range("B13").Value = "94641241A"

# Move the active selection to B14 (reflects the selection change in the diff)
$ws.Range("B14").Select()
